$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3C")

# Fix C3: was stored as a numeric-looking inline string "39", should become a genuine number 39
$ws.Range("C3").Value = 39

# Append the new submission as row 4
$ws.Range("A4").Value = "2026-02-07 20:00:19"
$ws.Range("B4").Value = "Muhammad Ali zarami"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "23"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = 8
